# Generate Report for Handoff
# Update "Latest Handoff Date/Datetime" for files that were just handed off
# (status: "Handback transform failed" and "Ready for handoff") across the
# Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# Rows 4, 6, 7, 8, 9, 10 correspond to:
#   bc392b61-e2ef-48ed-a1bd-14f89c7447d2.md
#   2b1022b4-be64-433c-9eb4-05221c97e060.md
#   4863f28a-9a37-4b27-901c-72949c36301a.md
#   58a7ffd7-3a9a-4d05-8d70-8753e3a0e1d7.md
#   5bfb8e10-718b-4fca-ad82-8483f9cb5a83.md
#   92105bf3-00f2-4e07-a7e0-cded7c381c71.md
$rows = @(4, 6, 7, 8, 9, 10)

foreach ($r in $rows) {
    $ws1.Range("D$r").Value = "2016-03-24 03:22:02"
    $ws2.Range("E$r").Value = "2016-03-24 03:21:53"
    $ws3.Range("E$r").Value = "2016-03-24 03:22:02"
}
